# Add four new scenario rows (37-40) to the Parameters sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Labels for the two "uptake" scenario rows first.
$ws.Range("A37").Value = "highup"
$ws.Range("A38").Value = "lowup"

# Then the parameter names for all four new rows.
$ws.Range("B37").Value = "HighUptakePercent"
$ws.Range("B38").Value = "LowUptakePercent"
$ws.Range("B39").Value = "WithHighSub"
$ws.Range("B40").Value = "WithLowSub"

# Labels for the two "sub" scenario rows.
$ws.Range("A39").Value = "highsub"
$ws.Range("A40").Value = "lowsub"

# Numeric values for each row.
$ws.Range("F37").Value = 0.12
$ws.Range("K37").Value = 0.45
$ws.Range("P37").Value = 0.09
$ws.Range("U37").Value = 0.45

$ws.Range("F38").Value = 0.05
$ws.Range("K38").Value = 0.2
$ws.Range("P38").Value = 0.04
$ws.Range("U38").Value = 0.2

$ws.Range("F39").Value = 0.11
$ws.Range("K39").Value = 0.41
$ws.Range("P39").Value = 0.082
$ws.Range("U39").Value = 0.41

$ws.Range("F40").Value = 0.087
$ws.Range("K40").Value = 0.335
$ws.Range("P40").Value = 0.067
$ws.Range("U40").Value = 0.335

# Update frozen-pane top-left cell and selection to reflect the new rows.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("U41").Select()
